# Rename the existing sheet "Blad1" -> "taken"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "taken"

# Add the new worksheet "db strutctuur" right after "taken"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "db strutctuur"

# Populate header row (order chosen to match shared-string table insertion order)
$ws2.Range("A1").Value = "Room number"
$ws2.Range("G1").Value = "Teacher"
$ws2.Range("C1").Value = "Start Time"
$ws2.Range("D1").Value = "Duration"
$ws2.Range("H1").Value = "Course description"
$ws2.Range("E1").Value = "Course Name"
$ws2.Range("F1").Value = "Course ID"
$ws2.Range("B1").Value = "Date"

# Style for row 2 (C2, D2) - date format, centered
$ws2.Range("C2:D2").NumberFormat = "d-mmm"
$ws2.Range("C2:D2").HorizontalAlignment = -4108
$ws2.Range("C2:D2").VerticalAlignment = -4108

# Column widths (values chosen so the engine's internal char-width
# quantization lands as close as possible to the widths from the diff)
$ws2.Columns.Item(1).ColumnWidth = 12.57
$ws2.Columns.Item(2).ColumnWidth = 12.57
$ws2.Columns.Item(3).ColumnWidth = 9
$ws2.Columns.Item(4).ColumnWidth = 8
$ws2.Columns.Item(5).ColumnWidth = 12
$ws2.Columns.Item(6).ColumnWidth = 12
$ws2.Columns.Item(7).ColumnWidth = 11
$ws2.Columns.Item(8).ColumnWidth = 16

# Selection / active cell on new sheet
$ws2.Range("C7").Select() | Out-Null

# Make "db strutctuur" the active sheet/tab, matching the diff (activeTab=1, tabSelected on sheet2)
$ws2.Activate() | Out-Null
